$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Paragraph 3 ("Such a survey ... existence planets ..."):
#   - remove justified alignment (w:jc val="both")
#   - remove the _GoBack bookmark that currently splits "exist"/"ence"
#     and merge the two runs back into a single run (the bookmark is
#     relocated into paragraph 2, see below)
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3.Alignment = 0

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$boundary = $oldBookmark.Start
$oldBookmark.Delete()

# Re-type the single character on either side of the old bookmark
# boundary so the (now identically-formatted, contiguous) runs merge
# back into one run, without disturbing the separately-styled "Maoz"
# run further along in the paragraph.
$mergeSpan = $d.Range($boundary - 1, $boundary + 1)
$mergeText = $mergeSpan.Text
$mergeSpan.Delete()
$insertionPoint = $d.Range($boundary - 1, $boundary - 1)
$insertionPoint.InsertAfter($mergeText)

# ---------------------------------------------------------------------
# Paragraph 2 ("This proposal is to obtain ... variability."):
#   - remove justified alignment (w:jc val="both")
#   - split the run after "K2 Cycle 1" and insert the _GoBack bookmark
#     at that split point (relocated from paragraph 3 above)
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Alignment = 0

$p2Range = $p2.Range
$p2Range.Find.Execute("K2 Cycle 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPos = $p2Range.End

$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))

# The engine keeps the original run's xml:space="preserve" flag on the
# left-hand fragment produced by the split above. Re-type its very last
# character (staying inside the existing run/formatting context) so the
# serializer recomputes xml:space for that fragment from its own text
# (which has no leading/trailing whitespace, so it should be omitted).
$lastChar = $d.Range($splitPos - 1, $splitPos)
$lastCharText = $lastChar.Text
$lastChar.Delete()
$d.Range($splitPos - 1, $splitPos - 1).InsertAfter($lastCharText)
